# gestionnaire emploi probleme de cache du controlleur
# A new "Affectation" row was missing from the schedule export because of a
# stale controller cache; re-insert it as row 3 (module 1DD 101 / Anglais
# technique / GEOCF / 2A / ESSADIK), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row just below the header/first data row, shifting the
# existing rows 3..35 down to 4..36.
$ws.Rows.Item(3).Insert()

# Fill in the newly inserted row with the missing affectation.
$ws.Cells.Item(3, 1).Value = "2024/2025"
$ws.Cells.Item(3, 2).Value = "1DD 101"
$ws.Cells.Item(3, 3).Value = "Anglais technique"
$ws.Cells.Item(3, 4).Value = "GEOCF"
$ws.Cells.Item(3, 5).Value = "2A"
$ws.Cells.Item(3, 6).Value = "ESSADIK"
$ws.Cells.Item(3, 7).Value = "ESSADIK"

# Match the reviewer's view state: zoomed to 85% with E3 as the active cell.
$excel.ActiveWindow.Zoom = 85
$ws.Range("E3").Select() | Out-Null
